$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are purely numeric-looking text (e.g. "1.00", "510.55")
# must be forced to Text format first, so Excel keeps them as literal strings
# instead of silently converting them to numbers (losing trailing zeros, etc.)
$textCells = @("D4","D5","D6","D7","D10","D13","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D33","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "58.083.99"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "2.458.62"
$ws.Range("E3").Value = "  +5.69%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "510.55"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "133.36"
$ws.Range("E6").Value = "  +7.65%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("D9").Value = "2.458.12"
$ws.Range("E9").Value = "  +8.12%  "
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  +5.41%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("D13").Value = "4.61"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "2.889.88"
$ws.Range("E14").Value = "  +4.59%  "
$ws.Range("D15").Value = "57.941.04"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "21.96"
$ws.Range("E16").Value = "  +5.76%  "
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  +6.91%  "
$ws.Range("D18").Value = "2.487.42"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "10.33"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").Value = "315.57"
$ws.Range("E20").Value = "  +5.21%  "
$ws.Range("D21").Value = "4.09"
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("D22").Value = "6.47"
$ws.Range("E22").Value = "  +10.32%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("D24").Value = "5.73"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").Value = "65.82"
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "0.156"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "0.382"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "7.57"
$ws.Range("E29").Value = "  +8.57%  "
$ws.Range("D30").Value = "171.18"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  +6.37%  "
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").Value = "6.10"
$ws.Range("E33").Value = "  +3.15%  "
$ws.Range("E34").Value = "  +4.38%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "18.11"
$ws.Range("E37").Value = "  +5.99%  "
$ws.Range("D38").Value = "1.22"
$ws.Range("E38").Value = "  +6.92%  "
$ws.Range("D39").Value = "3.89"
$ws.Range("E39").Value = "  +8.09%  "
$ws.Range("D40").Value = "36.70"
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  +7.36%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.47"
$ws.Range("E42").Value = "  +6.59%  "
$ws.Range("D43").Value = "135.82"
$ws.Range("E43").Value = "  +15.42%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.39"
$ws.Range("E44").Value = "  +5.36%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "4.92"
$ws.Range("E45").Value = "  +6.75%  "
$ws.Range("D46").Value = "255.67"
$ws.Range("E46").Value = "  +5.01%  "
$ws.Range("D47").Value = "0.573"
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("D49").Value = "0.0494"
$ws.Range("E49").Value = "  +4.04%  "
$ws.Range("D50").Value = "0.0213"
$ws.Range("E50").Value = "  +5.81%  "
$ws.Range("D51").Value = "17.24"
$ws.Range("E51").Value = "  +6.31%  "
